$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(3, 8, 4, 12),
    @(5, 16, 4, 4),
    @(3, 19, 4, 1),
    @(6, 7, 4, 13),
    @(4, 7, 1, 13),
    @(4, 13, 5, 7),
    @(5, 12, 7, 8),
    @(5, 12, 3, 8),
    @(4, 17, 2, 3),
    @(7, 16, 6, 4),
    @(5, 7, 6, 13),
    @(4, 8, 1, 12),
    @(4, 6, 3, 14),
    @(3, 15, 2, 5),
    @(4, 2, 3, 18),
    @(6, 4, 4, 16),
    @(9, 15, 5, 5),
    @(3, 12, 4, 8),
    @(2, 5, 4, 15)
)

$startRow = 247
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Formula = "=B$row+D$row"
}

$lastRow = $startRow + $data.Count - 1
$nextRow = $lastRow + 1
$ws.Range("A$nextRow").Select()
try { $excel.ActiveWindow.ScrollRow = $startRow } catch {}
